$d = $word.ActiveDocument

# Replace the full "回复邮件地址：2624001227@qq.com" text with the shorter
# descriptive text "电子邮件地址为QQ邮箱".
$d.Content.Find.Execute("回复邮件地址：2624001227@qq.com", $true, $false, $false, $false, $false,
                         $true, 1, $false, "电子邮件地址为QQ邮箱", 2)

# The _GoBack bookmark sits right after that run, at the end of the
# paragraph. Remember its position, then temporarily remove it so we can
# append the e-mail address as a new run at the very end of the paragraph
# (after the bookmark's original spot) and finally restore the bookmark in
# its original place, ahead of the newly inserted run.
$bm = $d.Bookmarks.Item("_GoBack")
$bmPos = $bm.Start
$bm.Delete()

# Find the paragraph that used to hold the bookmark (it sits right before
# the paragraph mark) and locate the very end of its text.
$p = $d.Range($bmPos, $bmPos)
[void]$p.Expand(4)
$r = $d.Range($p.End - 1, $p.End - 1)
$r.InsertAfter("：2624001227@qq.com")

# Touch the formatting so the newly inserted text is emitted as its own
# <w:r> run (matching the source run's appearance) instead of being merged
# back into the preceding run.
$r.Font.Bold = $true
$r.Font.Bold = $false

# Restore the _GoBack bookmark at its original (pre-insert) position, which
# now sits right before the new run we just appended.
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
